# Clarify that the report-creation timestamp shown in cell E1 is in UTC.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "Date & Time Report Created (UTC)"
